$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.873.60'
$ws.Cells.Item(2, 5).Value = '  -0.32%  '

$ws.Cells.Item(3, 4).Value = '1.633.94'
$ws.Cells.Item(3, 5).Value = '  -0.33%  '

$ws.Cells.Item(4, 5).Value = '  -0.25%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '213.81'
$ws.Cells.Item(5, 5).Value = '  -0.58%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.5051'
$ws.Cells.Item(6, 5).Value = '  -0.41%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.002'
$ws.Cells.Item(7, 5).Value = '  -0.15%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2574'
$ws.Cells.Item(8, 5).Value = '  +0.40%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06342'
$ws.Cells.Item(9, 5).Value = '  -0.40%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '19.51'
$ws.Cells.Item(10, 5).Value = '  +0.13%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07732'
$ws.Cells.Item(11, 5).Value = '  -0.40%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '4.260'
$ws.Cells.Item(12, 5).Value = '  -0.48%  '

$ws.Cells.Item(13, 4).Value = '1.623.19'
$ws.Cells.Item(13, 5).Value = '  -1.43%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.5421'
$ws.Cells.Item(14, 5).Value = '  -0.52%  '

$ws.Cells.Item(15, 4).Value = '0.0₅7697'
$ws.Cells.Item(15, 5).Value = '  -1.56%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '63.99'
$ws.Cells.Item(16, 5).Value = '  -0.59%  '

$ws.Cells.Item(17, 4).Value = '25.874.25'
$ws.Cells.Item(17, 5).Value = '  -0.61%  '

$ws.Cells.Item(18, 5).Value = '  +0.00%  '

$ws.Cells.Item(19, 2).Value = 'Uniswap'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.415'
$ws.Cells.Item(19, 5).Value = '  -0.95%  '

$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '194.27'
$ws.Cells.Item(20, 5).Value = '  -1.79%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '9.872'
$ws.Cells.Item(21, 5).Value = '  -0.75%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.069'
$ws.Cells.Item(22, 5).Value = '  +0.36%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.004'
$ws.Cells.Item(23, 5).Value = '  -0.04%  '

$ws.Cells.Item(24, 5).Value = '  +1.15%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '142.57'
$ws.Cells.Item(25, 5).Value = '  +0.98%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.1238'
$ws.Cells.Item(26, 5).Value = '  +6.03%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '6.802'
$ws.Cells.Item(27, 5).Value = '  -1.16%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '15.57'
$ws.Cells.Item(28, 5).Value = '  -0.93%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.235'
$ws.Cells.Item(29, 5).Value = '  -0.24%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.04846'
$ws.Cells.Item(30, 5).Value = '  -3.15%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '3.224'
$ws.Cells.Item(31, 5).Value = '  -0.88%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.177'
$ws.Cells.Item(32, 5).Value = '  -0.37%  '

$ws.Cells.Item(33, 5).Value = '  -0.03%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.367'
$ws.Cells.Item(34, 5).Value = '  +0.28%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.9033'
$ws.Cells.Item(35, 5).Value = '  +0.66%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.572'
$ws.Cells.Item(36, 5).Value = '  -0.55%  '

$ws.Cells.Item(37, 4).Value = '1.124.24'
$ws.Cells.Item(37, 5).Value = '  -0.20%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.5470'
$ws.Cells.Item(38, 5).Value = '  -0.05%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.01550'
$ws.Cells.Item(39, 5).Value = '  -0.38%  '

$ws.Cells.Item(40, 5).Value = '  -0.08%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '5.574'
$ws.Cells.Item(41, 5).Value = '  -0.66%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.8011'
$ws.Cells.Item(42, 5).Value = '  -2.10%  '

$ws.Cells.Item(43, 2).Value = 'Quant'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '98.29'
$ws.Cells.Item(43, 5).Value = '  -1.56%  '

$ws.Cells.Item(44, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(44, 4).Value = '0.0₈121'
$ws.Cells.Item(44, 5).Value = '  -6.69%  '

$ws.Cells.Item(45, 4).Value = '1.769.77'
$ws.Cells.Item(45, 5).Value = '  -0.79%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.4466'
$ws.Cells.Item(46, 5).Value = '  -1.54%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.001'
$ws.Cells.Item(47, 5).Value = '  -0.52%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '54.77'
$ws.Cells.Item(48, 5).Value = '  +0.04%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.05154'
$ws.Cells.Item(49, 5).Value = '  +1.64%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '7.539'
$ws.Cells.Item(50, 5).Value = '  +2.09%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.001'
$ws.Cells.Item(51, 5).Value = '  -0.27%  '
